$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the computed (0.196 * previous-row) formula rows, leaving the
# styled-but-empty cells behind (matches the target OOXML exactly).
$rowsToClear = @(4, 7, 10, 13, 16, 19, 22, 25, 28, 31, 34, 37)
foreach ($r in $rowsToClear) {
    [void]$ws.Range("A$r`:P$r").ClearContents()
}

# Update the frozen-pane view state: scroll/top-left cell and the
# selected (active) cell within the frozen bottom-left pane.
[void]$ws.Range("B36").Select()
